$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and reorder FraxShare/EnergySwap/ApeXProtocol rows)
# Values are prefixed with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr cell type) rather than reinterpreting numeric-looking
# strings like "33.08" as actual numbers.

$ws.Range("D2").Value = "'42.950.23"
$ws.Range("E2").Value = "'  -0.42%  "
$ws.Range("D3").Value = "'2.302.44"
$ws.Range("E3").Value = "'  -0.09%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'306.15"
$ws.Range("E5").Value = "'  +1.89%  "
$ws.Range("D6").Value = "'97.59"
$ws.Range("E6").Value = "'  -0.47%  "
$ws.Range("E7").Value = "'  -1.82%  "
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("D9").Value = "'0.504"
$ws.Range("E9").Value = "'  -2.34%  "
$ws.Range("E10").Value = "'  -0.15%  "
$ws.Range("E11").Value = "'  +0.00%  "
$ws.Range("D12").Value = "'18.21"
$ws.Range("E12").Value = "'  +0.71%  "
$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = "'  +1.04%  "
$ws.Range("E14").Value = "'  -1.54%  "
$ws.Range("D15").Value = "'2.660.96"
$ws.Range("E15").Value = "'  -0.17%  "
$ws.Range("D16").Value = "'2.305.98"
$ws.Range("E16").Value = "'  -2.56%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "'  -0.27%  "
$ws.Range("D18").Value = "'42.875.84"
$ws.Range("E18").Value = "'  -0.34%  "
$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "'  -4.57%  "
$ws.Range("E20").Value = "'  -0.51%  "
$ws.Range("E21").Value = "'  -1.15%  "
$ws.Range("E22").Value = "'  -1.14%  "
$ws.Range("D23").Value = "'237.02"
$ws.Range("E23").Value = "'  -0.54%  "
$ws.Range("E24").Value = "'  -2.47%  "
$ws.Range("E25").Value = "'  +1.66%  "
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("E27").Value = "'  +0.10%  "
$ws.Range("D28").Value = "'25.49"
$ws.Range("E28").Value = "'  +2.71%  "
$ws.Range("D29").Value = "'166.47"
$ws.Range("E29").Value = "'  -0.96%  "
$ws.Range("E30").Value = "'  +0.88%  "
$ws.Range("D31").Value = "'9.07"
$ws.Range("E31").Value = "'  -1.01%  "
$ws.Range("D32").Value = "'33.08"
$ws.Range("E33").Value = "'  +0.12%  "
$ws.Range("D34").Value = "'4.84"
$ws.Range("E34").Value = "'  +1.24%  "
$ws.Range("D35").Value = "'5.01"
$ws.Range("E35").Value = "'  -3.03%  "
$ws.Range("D36").Value = "'17.21"
$ws.Range("E36").Value = "'  -4.92%  "
$ws.Range("D37").Value = "'2.40"
$ws.Range("E37").Value = "'  -0.38%  "
$ws.Range("E38").Value = "'  +0.74%  "
$ws.Range("E39").Value = "'  -0.98%  "
$ws.Range("D40").Value = "'1.75"
$ws.Range("E40").Value = "'  -1.96%  "
$ws.Range("E41").Value = "'  -1.21%  "
$ws.Range("E42").Value = "'  -0.75%  "
$ws.Range("D43").Value = "'2.006.16"
$ws.Range("E43").Value = "'  -0.06%  "
$ws.Range("E44").Value = "'  -2.22%  "
$ws.Range("B45").Value = "'ApeXProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.13"
$ws.Range("E45").Value = "'  -1.57%  "
$ws.Range("B46").Value = "'FraxShare"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.02"
$ws.Range("E46").Value = "'  -1.48%  "
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.99"
$ws.Range("E47").Value = "'  +3.96%  "
$ws.Range("E48").Value = "'  -1.60%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "'  +5.24%  "
$ws.Range("D50").Value = "'54.13"
$ws.Range("E50").Value = "'  -0.35%  "
$ws.Range("D51").Value = "'2.528.13"
$ws.Range("E51").Value = "'  -0.07%  "
